# update P50 class list
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New first-name-only roster (replaces the old "Last, First" + first-name pairs)
$names = @(
    "Teddy",
    "Christine",
    "Max",
    "Wallace",
    "Smith",
    "Harrison",
    "Catherine",
    "Katie",
    "Jackson",
    "Izaac",
    "Jaden",
    "Christine",
    "Chris",
    "Zhanzhiz",
    "Ryan"
)

for ($i = 0; $i -lt $names.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $null
}

# Give the name column (A2:A16) its own smaller font, distinct from the rest of the roster rows
$ws.Range("A2").Font.Name = "Arial"
$ws.Range("A2").Font.Size = 10
$ws.Range("A2").Copy()
$ws.Range("A3:A16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.PageSetup.Orientation = 1

[void]$ws.Range("A17").Select()
